$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "Company Name"
$ws.Range("B1").Value = "Company Address"
$ws.Range("C1").Value = "Company Town"
$ws.Range("D1").Value = "Company State"
$ws.Range("E1").Value = "Company Contact"
$ws.Range("F1").Value = "Company TIN"
$ws.Range("G1").Value = "Company PIN"
$ws.Range("H1").Value = "Customer Name"
$ws.Range("I1").Value = "Customer Address"
$ws.Range("J1").Value = "Customer Town"
$ws.Range("K1").Value = "Customer State"
$ws.Range("L1").Value = "Customer PIN"
$ws.Range("M1").Value = "Customer TIN"
$ws.Range("N1").Value = "Customer Contact"
$ws.Range("O1").Value = "Invoice Number"
$ws.Range("P1").Value = "Invoice Date"
$ws.Range("Q1").Value = "Due Date"
$ws.Range("R1").Value = "Sub Total"
$ws.Range("S1").Value = "GST 8%"
$ws.Range("T1").Value = "Total"

# Row 2 - values
$ws.Range("A2").Value = "J.K. Computers"
$ws.Range("B2").Value = "1133 South Cavalier Drive"
$ws.Range("C2").Value = "Alamo"
$ws.Range("D2").Value = " Tennessee"

# These values look numeric, so the target cells must be forced to Text
# format first - otherwise Excel would silently convert them to numbers
# and drop the leading "+", leading/trailing spaces, etc.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "+17319696651"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "14123456789"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "343567"

$ws.Range("H2").Value = "Giana Kousky"
$ws.Range("I2").Value = "Roseville Building Room No.102"
$ws.Range("J2").Value = "Bellmont"
$ws.Range("K2").Value = " Illinois"

$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "5646"

$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "1412345678"

$ws.Range("N2").NumberFormat = "@"
$ws.Range("N2").Value = "+3435677789"

$ws.Range("O2").NumberFormat = "@"
$ws.Range("O2").Value = " 100"

$ws.Range("P2").Value = "Feb 23- 2016"
$ws.Range("Q2").Value = "Mar 10- 2016"

$ws.Range("R2").NumberFormat = "@"
$ws.Range("R2").Value = " 10800.00"

$ws.Range("S2").NumberFormat = "@"
$ws.Range("S2").Value = "864.00"

$ws.Range("T2").NumberFormat = "@"
$ws.Range("T2").Value = "11664.00 "
